$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.219.45'
$ws.Range('E2').Value = '  +0.23%  '
$ws.Range('D3').Value = '1.862.47'
$ws.Range('E3').Value = '  +0.76%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7039'
$ws.Range('E5').Value = '  +0.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '237.74'
$ws.Range('E6').Value = '  -0.28%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9993'
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.08221'
$ws.Range('E8').Value = '  +10.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3037'
$ws.Range('E9').Value = '  -0.54%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.30'
$ws.Range('E10').Value = '  -0.35%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08166'
$ws.Range('E11').Value = '  +0.51%  '
$ws.Range('D12').Value = '1.830.29'
$ws.Range('E12').Value = '  -1.78%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.165'
$ws.Range('E13').Value = '  -1.02%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.7094'
$ws.Range('E14').Value = '  -2.19%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '89.25'
$ws.Range('E15').Value = '  +0.66%  '
$ws.Range('D16').Value = '29.220.68'
$ws.Range('E16').Value = '  -0.29%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000007900'
$ws.Range('E17').Value = '  +3.65%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.784'
$ws.Range('E18').Value = '  +0.41%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.35'
$ws.Range('E19').Value = '  +2.26%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '236.41'
$ws.Range('E20').Value = '  -0.76%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9988'
$ws.Range('E21').Value = '  -0.14%  '
$ws.Range('D22').Value = '2.106.18'
$ws.Range('E22').Value = '  -0.94%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9993'
$ws.Range('E23').Value = '  -0.08%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.405'
$ws.Range('E24').Value = '  -2.31%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '162.18'
$ws.Range('E25').Value = '  +0.63%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.956'
$ws.Range('E26').Value = '  -0.39%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1438'
$ws.Range('E27').Value = '  -0.94%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.09'
$ws.Range('E28').Value = '  +0.12%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.961'
$ws.Range('E29').Value = '  -1.15%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.424'
$ws.Range('E30').Value = '  +1.85%  '
$ws.Range('E31').Value = '  -0.50%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.391'
$ws.Range('E32').Value = '  -3.36%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.056'
$ws.Range('E33').Value = '  +2.13%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05204'
$ws.Range('E34').Value = '  +0.60%  '
$ws.Range('E35').Value = '  -1.43%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7067'
$ws.Range('E36').Value = '  +0.94%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9993'
$ws.Range('E37').Value = '  -4.10%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.668'
$ws.Range('E38').Value = '  +0.50%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01851'
$ws.Range('E39').Value = '  -0.66%  '
$ws.Range('E40').Value = '  +1.94%  '
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').Value = '1.144.44'
$ws.Range('E41').Value = '  +6.26%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9240'
$ws.Range('E42').Value = '  -0.96%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4273'
$ws.Range('E43').Value = '  -0.22%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.868'
$ws.Range('E44').Value = '  -2.34%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '70.15'
$ws.Range('E45').Value = '  -0.19%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.9985'
$ws.Range('E46').Value = '  -0.13%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '102.80'
$ws.Range('E47').Value = '  +0.22%  '
$ws.Range('E48').Value = '  +1.87%  '
$ws.Range('D49').Value = '1.999.77'
$ws.Range('E49').Value = '  -0.54%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.199'
$ws.Range('E50').Value = '  +0.58%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.957'
$ws.Range('E51').Value = '  -0.99%  '
